$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data (columns D and E) for rows 2-51
# D-column values are forced to remain plain text (matching original inlineStr cells)
# by temporarily switching the cell to Text number format and restoring the original
# style afterwards, so Excel does not auto-convert strings like "1.002" into numbers.

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.924.38"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  -0.96%  "

$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.813.89"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -0.87%  "

$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  +0.14%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.38"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -2.07%  "

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5913"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -2.94%  "

$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  +0.19%  "

$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2749"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  -2.40%  "

$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06746"
$ws.Range("D9").Style = $style

$ws.Range("E10").Value = "  -3.91%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07488"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -2.01%  "

$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.819.45"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -0.54%  "

$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.673"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  -2.90%  "

$ws.Range("E14").Value = "  -1.62%  "

$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009282"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -7.01%  "

$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "74.48"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -6.17%  "

$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "28.655.61"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -1.91%  "

$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.441"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -8.56%  "

$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +0.14%  "

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.25"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  -8.90%  "

$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.37"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -3.71%  "

$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.772"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -3.80%  "

$ws.Range("E23").Value = "  +0.40%  "

$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "154.57"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -0.58%  "

$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1274"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -2.22%  "

$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.794"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -3.77%  "

$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.29"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -2.51%  "

$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06309"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -6.78%  "

$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.406"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -4.99%  "

$ws.Range("E30").Value = "  -2.10%  "

$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.729"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -3.16%  "

$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.684"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -4.04%  "

$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.690"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -2.74%  "

$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.050"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -7.10%  "

$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6348"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -3.12%  "

$ws.Range("E36").Value = "  -1.48%  "

$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.729"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -1.20%  "

$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.427"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -2.40%  "

$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01688"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -4.35%  "

$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.130.74"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -8.42%  "

$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8683"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -5.93%  "

$ws.Range("E42").Value = "  +0.15%  "

$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.968.59"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -0.91%  "

$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.02"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -0.93%  "

$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "60.46"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -4.96%  "

$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000112"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -4.16%  "

$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.571"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -3.44%  "

$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4511"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -1.16%  "

$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05452"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -1.94%  "

$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.277"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -3.60%  "

$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.003"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +0.21%  "
